$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: delete trailing rows 22 and 23 (sheet shrinks from 23 to 21 rows) ----
$ws.Rows.Item(22).Delete() | Out-Null
$ws.Rows.Item(22).Delete() | Out-Null

# ---- Step 2: clear cells that become empty (B17/C17) ----
$ws.Range("B17:C17").ClearContents() | Out-Null

# ---- Step 3: update cell values to match the target content/layout ----
$ws.Range("B10").Value = "5817330 - Larissa de Freitas"
$ws.Range("C10").Value = "5817330 - Larissa de Freitas"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1 - Introduction to the Chemistry Laboratory; 2 - Weights and measures; 3 - Methods for separating mixtures; 4 - Physical phenomena;  5 - Miscibility; 6 - Chemical Reactions; 7 -  Solutions; 8- Titrimetry; 9 - Chemical Equilibrium."
$ws.Range("C14").Value = "1 - Introduction to the Chemistry Laboratory; 2 - Weights and measures; 3 - Methods for separating mixtures; 4 - Physical phenomena;  5 - Miscibility; 6 - Chemical Reactions; 7 -  Solutions; 8- Titrimetry; 9 - Chemical Equilibrium."
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5817330 - Larissa de Freitas"
$ws.Range("C15").Value = "5817330 - Larissa de Freitas"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment.2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation.3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation.4 - Physical phenomena: Water state changes.5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects.7 - Solutions: Preparation and standardization of solutions.8 - Titrimetry: Acid-Base Titrations and return-titration.9 - Chemical equilibrium: Buffer solution."
$ws.Range("C16").Value = "1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment.2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation.3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation.4 - Physical phenomena: Water state changes.5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects.7 - Solutions: Preparation and standardization of solutions.8 - Titrimetry: Acid-Base Titrations and return-titration.9 - Chemical equilibrium: Buffer solution."
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "6310296 - Patrícia Caroline Molgero Da Rós"
$ws.Range("C18").Value = "6310296 - Patrícia Caroline Molgero Da Rós"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final será calculada da seguinte forma: NF = (3xMR + 7xMP)/10 onde NF é a nota final , MR é a média dos relatórios e MP é a média simples das provas."
$ws.Range("C20").Value = "A nota final será calculada da seguinte forma: NF = (3xMR + 7xMP)/10 onde NF é a nota final , MR é a média dos relatórios e MP é a média simples das provas."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova (PR) para alunos que tenham NF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a nota final (NF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova (PR) para alunos que tenham NF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a nota final (NF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0"

# ---- Step 3b: B13/C13 need "01/01/2018" kept as TEXT, not auto-converted to a date ----
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2018"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2018"
# restore original (General, wrap-top) cell formatting/style now that the text is already stored
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 3c: B19/C19 are brand-new cells; set values then fix up style/format to match column B/C ----
$ws.Range("B19").Value = "Os instrumentos de avaliação utilizados serão duas provas (P1 e P2) e a média dos relatórios (MR). O professor poderá a seu critério utilizar de trabalhos e/ou testes para complementar o método avaliativo."
$ws.Range("C19").Value = "Os instrumentos de avaliação utilizados serão duas provas (P1 e P2) e a média dos relatórios (MR). O professor poderá a seu critério utilizar de trabalhos e/ou testes para complementar o método avaliativo."
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 4: fix row heights ----
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
